# Fruta / hortaliza, semanal
# A new weekly price record (Femacal de La Calera, Papaya, "Primera" quality,
# 2023-03-23) is inserted before the existing row 61, pushing the old rows
# 61-64 down to rows 62-65 and extending the used range to A1:T65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61 - shifts rows 61:64 down to 62:65.
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with the new record's data.
$ws.Cells.Item(61, 1).Value = 3
$ws.Cells.Item(61, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(61, 3).Value = "Coquimbo"
$ws.Cells.Item(61, 4).Value = 45008
$ws.Cells.Item(61, 5).Value = 5
$ws.Cells.Item(61, 6).Value = "Fruta"
$ws.Cells.Item(61, 7).Value = 100108
$ws.Cells.Item(61, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(61, 9).Value = 100108004
$ws.Cells.Item(61, 10).Value = "Papaya"
$ws.Cells.Item(61, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(61, 12).Value = "Primera"
$ws.Cells.Item(61, 13).Value = 45
$ws.Cells.Item(61, 14).Value = 20000
$ws.Cells.Item(61, 15).Value = 20000
$ws.Cells.Item(61, 16).Value = 20000
$ws.Cells.Item(61, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(61, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(61, 19).Value = 2000
$ws.Cells.Item(61, 20).Value = 10
